$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, [string]$value)
    $style = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $style
}

Set-CellText $ws.Range("D2") '41.795.13'
Set-CellText $ws.Range("E2") '  +0.65%  '
Set-CellText $ws.Range("D3") '2.478.04'
Set-CellText $ws.Range("E3") '  +0.32%  '
Set-CellText $ws.Range("E4") '  -0.01%  '
Set-CellText $ws.Range("D5") '319.43'
Set-CellText $ws.Range("E5") '  +1.53%  '
Set-CellText $ws.Range("D6") '93.31'
Set-CellText $ws.Range("E6") '  +1.49%  '
Set-CellText $ws.Range("D7") '0.553'
Set-CellText $ws.Range("E7") '  +0.51%  '
Set-CellText $ws.Range("E8") '  +0.09%  '
Set-CellText $ws.Range("D9") '0.519'
Set-CellText $ws.Range("E9") '  +0.53%  '
Set-CellText $ws.Range("D10") '0.0883'
Set-CellText $ws.Range("E10") '  +11.42%  '
Set-CellText $ws.Range("D11") '33.29'
Set-CellText $ws.Range("E11") '  +2.39%  '
Set-CellText $ws.Range("E12") '  +0.75%  '
Set-CellText $ws.Range("D13") '2.860.74'
Set-CellText $ws.Range("E13") '  +0.33%  '
Set-CellText $ws.Range("E14") '  +1.09%  '
Set-CellText $ws.Range("D15") '15.76'
Set-CellText $ws.Range("E15") '  -1.55%  '
Set-CellText $ws.Range("D16") '2.490.45'
Set-CellText $ws.Range("E16") '  +0.79%  '
Set-CellText $ws.Range("E17") '  +2.80%  '
Set-CellText $ws.Range("D18") '41.753.68'
Set-CellText $ws.Range("B19") 'ShibaInu'
Set-CellText $ws.Range("C19") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText $ws.Range("D19") '0.0₃0954'
Set-CellText $ws.Range("E19") '  +1.22%  '
Set-CellText $ws.Range("B20") 'Uniswap'
Set-CellText $ws.Range("C20") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText $ws.Range("D20") '6.46'
Set-CellText $ws.Range("E20") '  -0.62%  '
Set-CellText $ws.Range("D21") '71.21'
Set-CellText $ws.Range("E21") '  +0.41%  '
Set-CellText $ws.Range("E22") '  +2.26%  '
Set-CellText $ws.Range("D23") '240.11'
Set-CellText $ws.Range("E23") '  +1.04%  '
Set-CellText $ws.Range("D24") '2.78'
Set-CellText $ws.Range("E24") '  +1.77%  '
Set-CellText $ws.Range("D25") '1.96'
Set-CellText $ws.Range("D27") '24.76'
Set-CellText $ws.Range("E27") '  +0.03%  '
Set-CellText $ws.Range("E28") '  +0.95%  '
Set-CellText $ws.Range("D29") '9.79'
Set-CellText $ws.Range("E29") '  +0.97%  '
Set-CellText $ws.Range("D30") '37.24'
Set-CellText $ws.Range("E30") '  +5.00%  '
Set-CellText $ws.Range("D31") '157.54'
Set-CellText $ws.Range("E31") '  +1.17%  '
Set-CellText $ws.Range("E32") '  +1.35%  '
Set-CellText $ws.Range("E33") '  +0.00%  '
Set-CellText $ws.Range("E34") '  +0.94%  '
Set-CellText $ws.Range("E35") '  +0.34%  '
Set-CellText $ws.Range("D36") '17.51'
Set-CellText $ws.Range("E36") '  +1.14%  '
Set-CellText $ws.Range("D37") '1.87'
Set-CellText $ws.Range("E37") '  +4.58%  '
Set-CellText $ws.Range("E38") '  +1.60%  '
Set-CellText $ws.Range("E39") '  +1.76%  '
Set-CellText $ws.Range("E40") '  +0.63%  '
Set-CellText $ws.Range("D41") '2.54'
Set-CellText $ws.Range("E41") '  +8.23%  '
Set-CellText $ws.Range("D42") '4.02'
Set-CellText $ws.Range("E42") '  +0.63%  '
Set-CellText $ws.Range("D43") '2.003.96'
Set-CellText $ws.Range("E43") '  +2.95%  '
Set-CellText $ws.Range("D44") '19.10'
Set-CellText $ws.Range("E44") '  +1.34%  '
Set-CellText $ws.Range("D45") '0.0285'
Set-CellText $ws.Range("E45") '  +0.70%  '
Set-CellText $ws.Range("E46") '  +3.21%  '
Set-CellText $ws.Range("D47") '9.52'
Set-CellText $ws.Range("E47") '  +4.71%  '
Set-CellText $ws.Range("D48") '2.717.00'
Set-CellText $ws.Range("E48") '  +0.29%  '
Set-CellText $ws.Range("D49") '98.57'
Set-CellText $ws.Range("E49") '  +1.26%  '
Set-CellText $ws.Range("D50") '74.65'
Set-CellText $ws.Range("E50") '  +4.25%  '
Set-CellText $ws.Range("D51") '67.41'
Set-CellText $ws.Range("E51") '  +0.33%  '
